# Taskmanagement.xlsx -- "Neues Projekt angelegt wegen visualisierung"
#
# The task list gained a new Task row ("Rausfinden wie Text rotieren") that
# was inserted above the existing "Datensimulator" / "Festlegen der
# Datentypen" entry, and a Status ("D") / Notiz ("E") column got filled in
# for a number of the existing tasks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fill in the Status column (D) for several tasks -------------------
$ws.Range("D2").Value = "Erledigt"
$ws.Range("D3").Value = "Erledigt"
$ws.Range("D4").Value = "Erledigt"
$ws.Range("D6").Value = "in Arbeit"
$ws.Range("D7").Value = "in Arbeit"

# --- 2. Add a note in column E (before the row insert below, so this new
# shared string lands right after the existing ones, matching the order
# the strings were authored in) --------------------------------------------
$ws.Range("E4").Value = "siehe Datenlexikon"

# Give the new Notiz column a sensible width, mirroring the bestFit columns
# already on the sheet.
$ws.Range("E1").EntireColumn.AutoFit()

# --- 3. Insert a new row at 14 -------------------------------------------
# This shifts every row from 14 downward by one (old row14 -> row15, ...,
# old row18 -> row19) and keeps all existing formatting/validation ranges
# attached the way a native Excel "Insert Row" would.
$ws.Rows("14:14").Insert()

# New task text for the freshly inserted row.
$ws.Range("B14").Value = "Rausfinden wie Text rotieren"

# --- 4. Fill in the Status column (D) for the rows that moved down --------
$ws.Range("D15").Value = "Erledigt"
$ws.Range("D16").Value = "in Arbeit"
$ws.Range("D17").Value = "in Arbeit"

# --- 5. Re-apply the Status dropdown validation over the full column ------
# (the insert above only preserves existing sqref pieces; make sure the
# whole data column, including the newly used rows, drives off the same
# $K$2:$K$5 list)
try { $ws.Range("D1:D1048576").Validation.Delete() } catch {}
$ws.Range("D22:D1048576").Validation.Add(3, 1, 1, "=`$K`$2:`$K`$5")
$ws.Range("D2:D21").Validation.Add(3, 1, 1, "=`$K`$2:`$K`$5")

# --- 6. Restore the cursor / selection state -------------------------------
$ws.Range("C9").Select()
